$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update matière field
$ws.Range("B3").Value = "devoir1"

# Update student last/first names (rows 10-19)
$ws.Range("A10").Value = "Armstrong"
$ws.Range("B10").Value = "Laura"

$ws.Range("A11").Value = "Bryant"
$ws.Range("B11").Value = "Michael"

$ws.Range("A12").Value = "Butler"
$ws.Range("B12").Value = "Timothy"

$ws.Range("A13").Value = "Carroll"
$ws.Range("B13").Value = "Theresa"

$ws.Range("A14").Value = "Cross"
$ws.Range("B14").Value = "Lauren"

$ws.Range("A15").Value = "Green"
$ws.Range("B15").Value = "Jesse"

$ws.Range("A16").Value = "Hernandez"
$ws.Range("B16").Value = "David"

$ws.Range("A17").Value = "Jensen"
$ws.Range("B17").Value = "Matthew"

$ws.Range("A18").Value = "Kelly"
$ws.Range("B18").Value = "Tyler"

$ws.Range("A19").Value = "Schmitt"
$ws.Range("B19").Value = "William"
